# Bill of materials updated:
#  - Add "Unit price" (G) and "Total price" (H) columns for every part row
#  - Add a TOTAL / exchange-rate (TRM) / COP-total block in rows 68-70
#  - Fix a handful of Digikey part-number typos in column E
#  - Minor page-setup tweaks (margins, header/footer) and view selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix Digikey reference typos in column E
# ---------------------------------------------------------------------------
$ws.Range("E27").Value = "490-1054-1-ND"          #     FB1
$ws.Range("E42").Value = "311-1.8KGRCT-ND"        # >  R4, R5, R7, R59
$ws.Range("E47").Value = "311-1.00KHRCT-ND"       # >  R3, R17-R20, R55  (trim leading space)
$ws.Range("E50").Value = "311-100HRCT-ND"         #     R58             (trim leading space)

# ---------------------------------------------------------------------------
# 2. New header cells for the price columns
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Unit price"
$ws.Range("H1").Value = "Total price"

# ---------------------------------------------------------------------------
# 3. Per-row unit price (G) and total price formula (H)
# ---------------------------------------------------------------------------
$unitPrices = @{
    2  = 0.1;    3  = 0.29;   4  = 0.11;   5  = 0.1;    6  = 0.1;
    7  = 1.68;   8  = 0.1;    9  = 0.1;    10 = 0.1;    11 = 0.1;
    12 = 0.1;    13 = 0.1;    14 = 0.1;    15 = 0.1;    16 = 0.1;
    17 = 0.1;    18 = 0.1;    19 = 0.1;    20 = 0.34;   21 = 0.13;
    22 = 0.29;   23 = 0.5;    24 = 0.31;   25 = 0.23;   26 = 0.73;
    27 = 0.1;    28 = 0.36;   29 = 0.89;   30 = 0.89;   31 = 0.72;
    32 = 0.42;   33 = 0.6;    34 = 0.44;   35 = 0.13;   36 = 0.1;
    37 = 0.17;   38 = 1.74;   39 = 0.51;   40 = 0.1;    41 = 0.1;
    42 = 0.1;    43 = 0.1;    44 = 0.1;    45 = 0.1;    46 = 0.1;
    47 = 0.1;    48 = 0.1;    49 = 0.1;    50 = 0.1;    51 = 0.1;
    52 = 0.48;   53 = 0.81;   54 = 0.35;   55 = 10.65;  56 = 2.42;
    57 = 11.21;  58 = 16.52;  59 = 10.74;  60 = 0.5;    61 = 0.5;
    62 = 0.74;   63 = 3.3;    64 = 0.68;   65 = 0.73;   66 = 0.1;
    67 = 0.14;
}

for ($r = 2; $r -le 67; $r++) {
    $ws.Cells.Item($r, 7).Value = $unitPrices[$r]
    $ws.Cells.Item($r, 8).Formula = "=+G$r*F$r"
}

# Currency number format for the new price columns
$ws.Range("G2:H67").NumberFormat = "[$`$-240A]#,##0.00;[RED]([$`$-240A]#,##0.00)"

# ---------------------------------------------------------------------------
# 4. TOTAL / TRM / COP summary block (rows 68-70)
# ---------------------------------------------------------------------------
$ws.Range("G68").Value = "TOTAL"
$ws.Range("H68").Formula = "=SUM(H2:H67)"
$ws.Range("I68").Value = "USD"

$ws.Range("G69").Value = "TRM"
$ws.Range("H69").Value = 3288

$ws.Range("H70").Formula = "=+H68*H69"
$ws.Range("I70").Value = "COP"

$ws.Range("G68:H70").NumberFormat = "[$`$-240A]#,##0.00;[RED]([$`$-240A]#,##0.00)"
$ws.Range("G68:I70").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. Page setup / header-footer tweaks
# ---------------------------------------------------------------------------
$ws.PageSetup.TopMargin = 73.8
$ws.PageSetup.BottomMargin = 73.8
$ws.PageSetup.CenterHeader = "&A"
$ws.PageSetup.CenterFooter = "Página &P"

# ---------------------------------------------------------------------------
# 6. View state - keep gridlines on, move selection/top-left cell
# ---------------------------------------------------------------------------
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("B1").Select()
$ws.Range("E22").Select()
